$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.339718954694745
$ws.Range("C2").Value = 4.715058708800846
$ws.Range("D2").Value = 4.778158895845936
$ws.Range("F2").Value = 23.23220913150073
$ws.Range("G2").Value = 3.625157938841989
$ws.Range("I2").Value = 19.9244378371014
$ws.Range("K2").Value = 7.459753489657254
$ws.Range("M2").Value = 19.97233895529696
$ws.Range("N2").Value = 18.23930963375482
$ws.Range("O2").Value = 20.82901816582712

$ws.Range("B3").Value = 7.088982912507766
$ws.Range("C3").Value = 4.593031919265494
$ws.Range("D3").Value = 4.710141732314316
$ws.Range("F3").Value = 23.23899888251899
$ws.Range("G3").Value = 3.626736601674625
$ws.Range("I3").Value = 19.99254590037702
$ws.Range("K3").Value = 7.231855891570631
$ws.Range("M3").Value = 19.37435090634192
$ws.Range("N3").Value = 18.30203441647594
$ws.Range("O3").Value = 20.88060037760386

$ws.Range("B4").Value = 6.93145942094254
$ws.Range("C4").Value = 4.515568212319816
$ws.Range("D4").Value = 4.667188016133835
$ws.Range("F4").Value = 23.2490296597983
$ws.Range("G4").Value = 3.62775704211213
$ws.Range("I4").Value = 20.0379485376675
$ws.Range("K4").Value = 7.086621068083423
$ws.Range("M4").Value = 19.00575640732409
$ws.Range("N4").Value = 18.34223040595374
$ws.Range("O4").Value = 20.91647812397404

$ws.Range("B5").Value = 6.866471928604363
$ws.Range("C5").Value = 4.483389803048322
$ws.Range("D5").Value = 4.649394308489276
$ws.Range("F5").Value = 23.25459066298024
$ws.Range("G5").Value = 3.628185780750041
$ws.Range("I5").Value = 20.05735068175375
$ws.Range("K5").Value = 7.026154063975852
$ws.Range("M5").Value = 18.8554454971488
$ws.Range("N5").Value = 18.35903536310085
$ws.Range("O5").Value = 20.93215425380357

$ws.Range("B6").Value = 6.855635937669399
$ws.Range("C6").Value = 4.478010512000288
$ws.Range("D6").Value = 4.646422469692236
$ws.Range("F6").Value = 23.25560301699474
$ws.Range("G6").Value = 3.628257752905625
$ws.Range("I6").Value = 20.06062673274559
$ws.Range("K6").Value = 7.016037692216731
$ws.Range("M6").Value = 18.83048774234955
$ws.Range("N6").Value = 18.36185151595545
$ws.Range("O6").Value = 20.93482096139893

$ws.Range("B7").Value = 6.930586055322411
$ws.Range("C7").Value = 4.51513668099134
$ws.Range("D7").Value = 4.666949203940452
$ws.Range("F7").Value = 23.24909869327492
$ws.Range("G7").Value = 3.627762771940081
$ws.Range("I7").Value = 20.03820655829905
$ws.Range("K7").Value = 7.08581071088858
$ws.Range("M7").Value = 19.00372934610894
$ws.Range("N7").Value = 18.34245532152088
$ws.Range("O7").Value = 20.91668526620762

$ws.Range("B8").Value = 7.25406488773421
$ws.Range("C8").Value = 4.673526395771961
$ws.Range("D8").Value = 4.754959486843142
$ws.Range("F8").Value = 23.23333347963972
$ws.Range("G8").Value = 3.625691674367992
$ws.Range("I8").Value = 19.94717702154939
$ws.Range("K8").Value = 7.382304655990171
$ws.Range("M8").Value = 19.7665969259161
$ws.Range("N8").Value = 18.26058887934326
$ws.Range("O8").Value = 20.84592958682119

$ws.Range("B9").Value = 7.85603473322446
$ws.Range("C9").Value = 4.962900089696745
$ws.Range("D9").Value = 4.91768910515078
$ws.Range("F9").Value = 23.24892274205633
$ws.Range("G9").Value = 3.622034047223297
$ws.Range("I9").Value = 19.79715163208351
$ws.Range("K9").Value = 7.919631377556419
$ws.Range("M9").Value = 21.2408181164338
$ws.Range("N9").Value = 18.11332470127791
$ws.Range("O9").Value = 20.74063559145747

$ws.Range("B10").Value = 8.273786535156253
$ws.Range("C10").Value = 5.161254345319088
$ws.Range("D10").Value = 5.030716853994358
$ws.Range("F10").Value = 23.28866781224117
$ws.Range("G10").Value = 3.61959022111111
$ws.Range("I10").Value = 19.70435437669943
$ws.Range("K10").Value = 8.285164122913709
$ws.Range("M10").Value = 22.29779543450976
$ws.Range("N10").Value = 18.01311439301302
$ws.Range("O10").Value = 20.68378141148171

$ws.Range("B11").Value = 8.45767001832388
$ws.Range("C11").Value = 5.248159553547413
$ws.Range("D11").Value = 5.080615867960413
$ws.Range("F11").Value = 23.31286097325441
$ws.Range("G11").Value = 3.618530738040646
$ws.Range("I11").Value = 19.66593577416839
$ws.Range("K11").Value = 8.444692431981149
$ws.Range("M11").Value = 22.77058408787295
$ws.Range("N11").Value = 17.96923647912542
$ws.Range("O11").Value = 20.66238965807381

$ws.Range("B12").Value = 8.52635765861449
$ws.Range("C12").Value = 5.280572455084932
$ws.Range("D12").Value = 5.09928496608867
$ws.Range("F12").Value = 23.32289709715491
$ws.Range("G12").Value = 3.618137004967607
$ws.Range("I12").Value = 19.65193454630322
$ws.Range("K12").Value = 8.504101024628881
$ws.Range("M12").Value = 22.94828634475358
$ws.Range("N12").Value = 17.95286492115574
$ws.Range("O12").Value = 20.65493350872615

$ws.Range("B13").Value = 8.51160745804761
$ws.Range("C13").Value = 5.273614090634261
$ws.Range("D13").Value = 5.095274446575558
$ws.Range("F13").Value = 23.32069681477848
$ws.Range("G13").Value = 3.618221470823491
$ws.Range("I13").Value = 19.65492560991379
$ws.Range("K13").Value = 8.491351295837781
$ws.Range("M13").Value = 22.91007688404012
$ws.Range("N13").Value = 17.95637999996776
$ws.Range("O13").Value = 20.65651064103329

$ws.Range("B14").Value = 8.463340250599689
$ws.Range("C14").Value = 5.250836222978581
$ws.Range("D14").Value = 5.082156372899343
$ws.Range("F14").Value = 23.31366914908941
$ws.Range("G14").Value = 3.618498195898359
$ws.Range("I14").Value = 19.66477291337851
$ws.Range("K14").Value = 8.449600240013856
$ws.Range("M14").Value = 22.78523141802902
$ws.Range("N14").Value = 17.96788469907543
$ws.Range("O14").Value = 20.66176331294975

$ws.Range("B15").Value = 8.433650459151343
$ws.Range("C15").Value = 5.23681900174686
$ws.Range("D15").Value = 5.074091427800031
$ws.Range("F15").Value = 23.30947826795825
$ws.Range("G15").Value = 3.618668669654567
$ws.Range("I15").Value = 19.67087595315589
$ws.Range("K15").Value = 8.423895280887978
$ws.Range("M15").Value = 22.70858147089843
$ws.Range("N15").Value = 17.9749633955031
$ws.Range("O15").Value = 20.66506469148869

$ws.Range("B16").Value = 8.261640178443496
$ws.Range("C16").Value = 5.155506375359437
$ws.Range("D16").Value = 5.027424572886624
$ws.Range("F16").Value = 23.28720939106404
$ws.Range("G16").Value = 3.619660508327134
$ws.Range("I16").Value = 19.7069416058711
$ws.Range("K16").Value = 8.274600007160139
$ws.Range("M16").Value = 22.26672035603697
$ws.Range("N16").Value = 18.01601615655191
$ws.Range("O16").Value = 20.68526951881478

$ws.Range("B17").Value = 8.154496161545358
$ws.Range("C17").Value = 5.104758343956981
$ws.Range("D17").Value = 4.998401161499554
$ws.Range("F17").Value = 23.2751110392194
$ws.Range("G17").Value = 3.620282317083106
$ws.Range("I17").Value = 19.73003971988568
$ws.Range("K17").Value = 8.181259708885042
$ws.Range("M17").Value = 21.99346326015997
$ws.Range("N17").Value = 18.0416370911264
$ws.Range("O17").Value = 20.69881087232133

$ws.Range("B18").Value = 8.092293630243169
$ws.Range("C18").Value = 5.075257200191944
$ws.Range("D18").Value = 4.981565263553479
$ws.Range("F18").Value = 23.26872808421253
$ws.Range("G18").Value = 3.620644883131887
$ws.Range("I18").Value = 19.74368230221713
$ws.Range("K18").Value = 8.126938642171133
$ws.Range("M18").Value = 21.83554612216476
$ws.Range("N18").Value = 18.05653445958252
$ws.Range("O18").Value = 20.70702030440664

$ws.Range("B19").Value = 8.071135923891285
$ws.Range("C19").Value = 5.065215566543866
$ws.Range("D19").Value = 4.975840720525965
$ws.Range("F19").Value = 23.266665918419
$ws.Range("G19").Value = 3.62076848767007
$ws.Range("I19").Value = 19.74836274858886
$ws.Range("K19").Value = 8.108438513181177
$ws.Range("M19").Value = 21.78195549438026
$ws.Range("N19").Value = 18.06160613300666
$ws.Range("O19").Value = 20.70987210330563

$ws.Range("B20").Value = 8.165961910614145
$ws.Range("C20").Value = 5.11019300627417
$ws.Range("D20").Value = 5.001505552480732
$ws.Range("F20").Value = 23.27633937539846
$ws.Range("G20").Value = 3.620215615765106
$ws.Range("I20").Value = 19.72754391227599
$ws.Range("K20").Value = 8.191261808838824
$ws.Range("M20").Value = 22.02263058179232
$ws.Range("N20").Value = 18.03889305947924
$ws.Range("O20").Value = 20.69732580947697

$ws.Range("B21").Value = 8.477543590972211
$ws.Range("C21").Value = 5.257540235009735
$ws.Range("D21").Value = 5.086015682437554
$ws.Range("F21").Value = 23.31570964495079
$ws.Range("G21").Value = 3.618416712614002
$ws.Range("I21").Value = 19.66186566422813
$ws.Range("K21").Value = 8.461890934463639
$ws.Range("M21").Value = 22.82193904743495
$ws.Range("N21").Value = 17.96449887914421
$ws.Range("O21").Value = 20.66020297526054

$ws.Range("B22").Value = 8.675648881505408
$ws.Range("C22").Value = 5.350940810055409
$ws.Range("D22").Value = 5.139923725648742
$ws.Range("F22").Value = 23.34653613690843
$ws.Range("G22").Value = 3.617284548950722
$ws.Range("I22").Value = 19.62213042739282
$ws.Range("K22").Value = 8.632914782948765
$ws.Range("M22").Value = 23.33649626911952
$ws.Range("N22").Value = 17.91729978576715
$ws.Range("O22").Value = 20.63969790474878

$ws.Range("B23").Value = 8.57044049969887
$ws.Range("C23").Value = 5.301361883950664
$ws.Range("D23").Value = 5.111275789669365
$ws.Range("F23").Value = 23.32961886999243
$ws.Range("G23").Value = 3.617884836703282
$ws.Range("I23").Value = 19.64304562487643
$ws.Range("K23").Value = 8.542180245384836
$ws.Range("M23").Value = 23.06263851299598
$ws.Range("N23").Value = 17.94236126085894
$ws.Range("O23").Value = 20.65029766524743

$ws.Range("B24").Value = 8.16078012217435
$ws.Range("C24").Value = 5.107737007698498
$ws.Range("D24").Value = 5.000102523751027
$ws.Range("F24").Value = 23.27578226093322
$ws.Range("G24").Value = 3.620245755605317
$ws.Range("I24").Value = 19.72867113580477
$ws.Range("K24").Value = 8.186741908169179
$ws.Range("M24").Value = 22.00944656992778
$ws.Range("N24").Value = 18.04013311418004
$ws.Range("O24").Value = 20.69799588481488

$ws.Range("B25").Value = 7.697179048485919
$ws.Range("C25").Value = 4.887032140277149
$ws.Range("D25").Value = 4.874776980054038
$ws.Range("F25").Value = 23.23972926682939
$ws.Range("G25").Value = 3.622980586668263
$ws.Range("I25").Value = 19.83468169009645
$ws.Range("K25").Value = 7.779251758285141
$ws.Range("M25").Value = 20.84572451798325
$ws.Range("N25").Value = 18.15175352684151
$ws.Range("O25").Value = 20.76552668531491
